$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("C2").Value = "350 CAMINO DE LA REINA STE 100  SAN DIEGO CA 92108 "
$ws.Range("E2").Value = "BUNCE MARTIN "
$ws.Range("F2").Value = "PO BOX 30968  MIDDLEBURG OH 44130 "

# G2/R2 contain digit-only text with a trailing space; prefix with a quote so
# Excel stores them as text (not auto-converted numbers), then reset the
# style so no stray NumberFormat/quote-prefix styling sticks to the cell.
$ws.Range("G2").Value = "'79862 "
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").Value = "ANTOINETTE GILKEY "
$ws.Range("I2").Value = "AKA ANTOINETTE A GILKEY 1929 WASHBURN ST CINCINNATI OH 45223  "

# J2 becomes a blank text cell (was "D 1 ")
$ws.Range("J2").Value = "'"
$ws.Range("J2").Style = "Normal"

$ws.Range("K2").Value = "23CV17210 "
$ws.Range("M2").Value = "MIDLAND CREDIT MANAGEMENT  INC vs. ANTOINETTE  GILKEY "

$ws.Range("R2").Value = "'6827.23 "
$ws.Range("R2").Style = "Normal"

# --- Row 3 ---
$ws.Range("C3").Value = "350 CAMINO DE LA REINA STE 100  SAN DIEGO CA 92108 "
$ws.Range("E3").Value = "PAVLOVIC NEVENKA "
$ws.Range("F3").Value = "PO BOX 2121  WARREN MI 48090 "

$ws.Range("G3").Value = "'72697 "
$ws.Range("G3").Style = "Normal"

$ws.Range("H3").Value = "KENARIYE DELANEY "
$ws.Range("I3").Value = "1410 SPRINGFIELD PIKE APT 2  CINCINNATI OH 45215 "
$ws.Range("K3").Value = "23CV17211 "
$ws.Range("M3").Value = "MIDLAND CREDIT MANAGEMENT  INC vs. KENARIYE  DELANEY "

$ws.Range("R3").Value = "'1867.21 "
$ws.Range("R3").Style = "Normal"
